# Fruta / hortaliza, semanal
#
# The weekly refresh prepends 3 new price observations (rows 373-375) to the
# "Frutilla" sheet, pushing the previously-existing rows 373-439 down to
# 376-442 (pure row-insert, no other content changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 373, shifting everything below
# (rows 373:439) down to 376:442.
$ws.Rows.Item(373).Resize(3).Insert()

# Values shared by every data row in this sheet (constant columns).
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId   = 100101
$producto     = "Berries"
$categoriaId  = 100112025
$categoria    = "Frutilla"
$variedad     = "Sin especificar"
$kgUnidad     = 7

# New rows data: Fecha, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm,
# UnidadComercializacion, Origen, Precio$/Kg
$newRows = @(
    @{ Row=373; Fecha=45209; Calidad="Especial"; Volumen=300; Min=13000; Max=13000; Prom=13000; Unidad="$/bandeja 7 kilos"; Origen="Provincia de Melipilla"; PrecioKg=1857 },
    @{ Row=374; Fecha=45209; Calidad="Especial"; Volumen=200; Min=11000; Max=11000; Prom=11000; Unidad="$/bandeja 7 kilos"; Origen="Provincia de Melipilla"; PrecioKg=1571 },
    @{ Row=375; Fecha=45209; Calidad="Segunda";  Volumen=200; Min=9000;  Max=9000;  Prom=9000;  Unidad="$/bandeja 7 kilos"; Origen="Provincia de Melipilla"; PrecioKg=1286 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
